$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# Slide 1 / Shape 1 (Title): "General thoughts and comments"
#   -> "Thoughts and comments", typed in-place as three runs:
#      "T" / "houghts " / "and comments"
# ---------------------------------------------------------------------
$titleTr = $s.Shapes.Item(1).TextFrame.TextRange

# Remove the leading "General " (first 8 characters).
$titleTr.Characters(1, 8).Text = ""

# Capitalise the now-leading "t" -> "T".
$titleTr.Characters(1, 1).Text = "T"

# Re-touch "houghts " so it becomes its own run, separate from the
# trailing "and comments" run.
$titleTr.Characters(2, 8).Text = "houghts "

# ---------------------------------------------------------------------
# Slide 1 / Shape 2 (Content placeholder): merge the sentence
# "This affects ... to that" and the trailing "." into a single run.
# (Text is set via a literal so the U+2019 right single quote in
# "don't" survives unchanged.)
# ---------------------------------------------------------------------
$bodyTr = $s.Shapes.Item(2).TextFrame.TextRange
$bodyFull = $bodyTr.Text
$needle = "This affects"
$idx = $bodyFull.IndexOf($needle)
$startPos = $idx + 1
$periodIdx = $bodyFull.IndexOf(".", $idx)
$len = $periodIdx - $idx + 1
$sentence = $bodyTr.Characters($startPos, $len)
$sentence.Text = "This affects the size of the dimension text when rescaled and I don’t know the solution to that."
